# Generate Report for Handback
# Update the timestamp values that get refreshed each time the handback
# status report is (re)generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row.
# (Shared with de-de!H2 which shows the same value.)
$overview.Range("G2").Value = "2016-08-18 02:58:05"

# de-de sheet mirrors the same "Latest HO Xliff Generate Date" value.
$dede.Range("H2").Value = "2016-08-18 02:58:05"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-18 02:57:57"
$zhcn.Range("K2").Value = "2016-08-18 02:58:26"

# de-de sheet: Correspond Handback DateTime
$dede.Range("K2").Value = "2016-08-18 02:58:33"
